$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a cell's value while forcing text interpretation
# (prevents Excel from auto-converting numeric-looking strings like
# "29.293.73" or "1.000" into actual numbers) and without altering
# the cell's existing style.
function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $origStyle = $r.Style
    $r.Value = "'" + $value
    $r.Style = $origStyle
}

# --- Simple price / volume(1h) updates (rows whose coin identity did not change) ---
$pairs = @(
    @("D2", "29.293.73"),
    @("E2", "  +0.56%  "),
    @("D3", "1.849.04"),
    @("E3", "  +0.38%  "),
    @("D4", "0.9988"),
    @("E4", "  -0.06%  "),
    @("D5", "241.30"),
    @("E5", "  +0.05%  "),
    @("D6", "0.6753"),
    @("E6", "  -1.66%  "),
    @("D7", "0.9997"),
    @("E7", "  -0.01%  "),
    @("D8", "0.07460"),
    @("E8", "  -0.12%  "),
    @("D9", "0.2961"),
    @("E9", "  -1.99%  "),
    @("D10", "23.04"),
    @("E10", "  -0.47%  "),
    @("D11", "0.07721"),
    @("E11", "  +0.69%  "),
    @("D12", "1.844.57"),
    @("E12", "  +0.15%  "),
    @("E13", "  -1.06%  "),
    @("D14", "0.6745"),
    @("E14", "  -1.31%  "),
    @("D15", "86.31"),
    @("E15", "  -1.38%  "),
    @("D16", "6.177"),
    @("E16", "  +0.06%  "),
    @("D17", "0.000008353"),
    @("E17", "  +2.33%  "),
    @("D18", "29.270.80"),
    @("E18", "  +0.56%  "),
    @("D19", "229.22"),
    @("E19", "  +0.36%  "),
    @("D21", "1.000"),
    @("E21", "  +0.04%  "),
    @("D22", "7.225"),
    @("E22", "  -2.41%  "),
    @("D23", "0.9995"),
    @("D24", "161.01"),
    @("E24", "  +0.63%  "),
    @("D25", "8.742"),
    @("E25", "  -0.24%  "),
    @("D26", "0.1413"),
    @("E26", "  -2.92%  "),
    @("E27", "  -0.24%  "),
    @("D28", "1.515"),
    @("E28", "  +0.38%  "),
    @("D29", "4.180"),
    @("E29", "  -1.98%  "),
    @("D30", "4.081"),
    @("E30", "  -1.41%  "),
    @("D31", "1.196"),
    @("E31", "  -0.03%  "),
    @("D32", "0.05328"),
    @("E32", "  +2.40%  "),
    @("D33", "0.7603"),
    @("E33", "  -0.75%  "),
    @("E34", "  +1.70%  "),
    @("D35", "1.141"),
    @("E35", "  +0.55%  "),
    @("D36", "2.674"),
    @("E36", "  -0.15%  "),
    @("D37", "1.323.46"),
    @("E37", "  +0.61%  "),
    @("E38", "  -1.47%  "),
    @("D39", "2.732"),
    @("E39", "  +0.14%  "),
    @("D40", "0.9202"),
    @("E40", "  -1.32%  "),
    @("E41", "  +3.59%  "),
    @("E42", "  +0.19%  "),
    @("D43", "103.67"),
    @("E43", "  -1.07%  "),
    @("D44", "0.08243"),
    @("E44", "  +11.52%  "),
    @("D45", "1.988.03"),
    @("E45", "  +0.27%  "),
    @("D50", "9.172"),
    @("E50", "  -3.60%  "),
    @("D51", "0.05958"),
    @("E51", "  +0.22%  ")

)

foreach ($p in $pairs) {
    Set-TextValue $p[0] $p[1]
}

# --- Rows 46-49 were reshuffled: BabyDogeCoin moved up to rank 44,
#     pushing Mantle / RenderToken / Aave each down one slot, with
#     refreshed price & volume figures for every row involved. ---
$rowUpdates = @(
    @{ Row = 46; Coin = "BabyDogeCoin"; Link = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"; Price = "0.00000000123"; Volume = "  +0.70%  " },
    @{ Row = 47; Coin = "Mantle";       Link = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt";              Price = "0.5169";        Volume = "  -0.60%  " },
    @{ Row = 48; Coin = "RenderToken";  Link = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr";    Price = "1.785";         Volume = "  +0.62%  " },
    @{ Row = 49; Coin = "Aave";         Link = "https://coinranking.com/coin/ixgUfzmLR+aave-aave";               Price = "64.28";         Volume = "  -1.09%  " }
)

foreach ($u in $rowUpdates) {
    Set-TextValue ("B" + $u.Row) $u.Coin
    Set-TextValue ("C" + $u.Row) $u.Link
    Set-TextValue ("D" + $u.Row) $u.Price
    Set-TextValue ("E" + $u.Row) $u.Volume
}
